# Update "想去人数" (want-to-go count) figures for the refreshed data pull.
$wb = $excel.ActiveWorkbook

# 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 410
$ws1.Range("F5").Value  = 869
$ws1.Range("F9").Value  = 13158
$ws1.Range("F10").Value = 1010
$ws1.Range("F15").Value = 664
$ws1.Range("F16").Value = 2087
$ws1.Range("F19").Value = 61
$ws1.Range("F21").Value = 232
$ws1.Range("F23").Value = 749

# 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 113

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 29

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 410
$ws4.Range("F6").Value  = 869
$ws4.Range("F11").Value = 13158
$ws4.Range("F12").Value = 1010
$ws4.Range("F17").Value = 664
$ws4.Range("F20").Value = 2087
$ws4.Range("F23").Value = 61
$ws4.Range("F27").Value = 29
$ws4.Range("F28").Value = 232
$ws4.Range("F30").Value = 749
$ws4.Range("F31").Value = 113
